$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 20 de Abril de 2020 a las 13:52'
$ws.Cells.Item(4, 2).Value = 56963
$ws.Cells.Item(4, 3).Value = 31762
$ws.Cells.Item(4, 4).Value = 17850
$ws.Cells.Item(4, 5).Value = 7351
$ws.Cells.Item(8, 2).Value = 6358
$ws.Cells.Item(8, 3).Value = 4178
$ws.Cells.Item(8, 4).Value = 10597
$ws.Cells.Item(8, 5).Value = 802
$ws.Cells.Item(11, 2).Value = 3938
$ws.Cells.Item(11, 3).Value = 4178
$ws.Cells.Item(11, 4).Value = 10597
$ws.Cells.Item(11, 5).Value = 504
$ws.Cells.Item(12, 2).Value = 3754
$ws.Cells.Item(12, 3).Value = 4178
$ws.Cells.Item(12, 4).Value = 10597
$ws.Cells.Item(12, 5).Value = 373
$ws.Cells.Item(13, 2).Value = 3734
$ws.Cells.Item(13, 3).Value = 1833
$ws.Cells.Item(13, 4).Value = 1616
$ws.Cells.Item(13, 5).Value = 285
$ws.Cells.Item(14, 2).Value = 3678
$ws.Cells.Item(14, 3).Value = 913
$ws.Cells.Item(14, 4).Value = 2274
$ws.Cells.Item(14, 5).Value = 491
$ws.Cells.Item(17, 2).Value = 3154
$ws.Cells.Item(17, 3).Value = 1070
$ws.Cells.Item(17, 4).Value = 1824
$ws.Cells.Item(17, 5).Value = 260
$ws.Cells.Item(19, 1).Value = 'Salamanca'
$ws.Cells.Item(19, 2).Value = 2602
$ws.Cells.Item(19, 3).Value = 794
$ws.Cells.Item(19, 4).Value = 1521
$ws.Cells.Item(19, 5).Value = 287
$ws.Cells.Item(20, 1).Value = 'Malaga'
$ws.Cells.Item(20, 2).Value = 2531
$ws.Cells.Item(20, 3).Value = 869
$ws.Cells.Item(20, 4).Value = 1439
$ws.Cells.Item(20, 5).Value = 223
$ws.Cells.Item(21, 1).Value = 'Segovia'
$ws.Cells.Item(21, 2).Value = 2406
$ws.Cells.Item(21, 3).Value = 656
$ws.Cells.Item(21, 4).Value = 1578
$ws.Cells.Item(21, 5).Value = 172
$ws.Cells.Item(22, 1).Value = 'Leon'
$ws.Cells.Item(22, 2).Value = 2403
$ws.Cells.Item(22, 3).Value = 1076
$ws.Cells.Item(22, 4).Value = 1024
$ws.Cells.Item(22, 5).Value = 303
$ws.Cells.Item(23, 1).Value = 'Asturias'
$ws.Cells.Item(23, 2).Value = 2348
$ws.Cells.Item(23, 3).Value = 599
$ws.Cells.Item(23, 4).Value = 1549
$ws.Cells.Item(23, 5).Value = 200
$ws.Cells.Item(24, 1).Value = 'Sevilla'
$ws.Cells.Item(24, 2).Value = 2329
$ws.Cells.Item(24, 3).Value = 459
$ws.Cells.Item(24, 4).Value = 1658
$ws.Cells.Item(24, 5).Value = 212
$ws.Cells.Item(25, 1).Value = 'Gipuzkoa/Guipuzcoa'
$ws.Cells.Item(25, 2).Value = 2328
$ws.Cells.Item(25, 3).Value = 6144
$ws.Cells.Item(25, 4).Value = 4953
$ws.Cells.Item(25, 5).Value = 209
$ws.Cells.Item(26, 2).Value = 2220
$ws.Cells.Item(26, 3).Value = 422
$ws.Cells.Item(26, 4).Value = 1482
$ws.Cells.Item(31, 1).Value = 'Burgos'
$ws.Cells.Item(31, 2).Value = 1567
$ws.Cells.Item(31, 3).Value = 642
$ws.Cells.Item(31, 4).Value = 757
$ws.Cells.Item(31, 5).Value = 168
$ws.Cells.Item(32, 1).Value = 'Pontevedra'
$ws.Cells.Item(32, 2).Value = 1536
$ws.Cells.Item(32, 3).Value = 333
$ws.Cells.Item(32, 4).Value = 1411
$ws.Cells.Item(32, 5).Value = 30
$ws.Cells.Item(33, 2).Value = 1431
$ws.Cells.Item(33, 3).Value = 4178
$ws.Cells.Item(33, 4).Value = 10597
$ws.Cells.Item(33, 5).Value = 186
$ws.Cells.Item(35, 1).Value = 'Cuenca'
$ws.Cells.Item(35, 2).Value = 1315
$ws.Cells.Item(35, 3).Value = 4178
$ws.Cells.Item(35, 4).Value = 10597
$ws.Cells.Item(35, 5).Value = 156
$ws.Cells.Item(36, 1).Value = 'Jaen'
$ws.Cells.Item(36, 2).Value = 1297
$ws.Cells.Item(36, 3).Value = 253
$ws.Cells.Item(36, 4).Value = 907
$ws.Cells.Item(36, 5).Value = 137
$ws.Cells.Item(39, 2).Value = 1243
$ws.Cells.Item(39, 3).Value = 299
$ws.Cells.Item(39, 4).Value = 848
$ws.Cells.Item(39, 5).Value = 96
$ws.Cells.Item(40, 1).Value = 'Avila'
$ws.Cells.Item(40, 2).Value = 1155
$ws.Cells.Item(40, 3).Value = 450
$ws.Cells.Item(40, 4).Value = 596
$ws.Cells.Item(40, 5).Value = 109
$ws.Cells.Item(41, 1).Value = 'Cadiz'
$ws.Cells.Item(41, 2).Value = 1139
$ws.Cells.Item(41, 3).Value = 283
$ws.Cells.Item(41, 4).Value = 781
$ws.Cells.Item(41, 5).Value = 75
$ws.Cells.Item(42, 2).Value = 1023
$ws.Cells.Item(42, 3).Value = 436
$ws.Cells.Item(42, 4).Value = 514
$ws.Cells.Item(42, 5).Value = 73
$ws.Cells.Item(45, 2).Value = 716
$ws.Cells.Item(45, 3).Value = 223
$ws.Cells.Item(45, 4).Value = 432
$ws.Cells.Item(45, 5).Value = 61
$ws.Cells.Item(46, 2).Value = 611
$ws.Cells.Item(46, 3).Value = 219
$ws.Cells.Item(46, 4).Value = 327
$ws.Cells.Item(46, 5).Value = 65
$ws.Cells.Item(47, 1).Value = 'Huesca'
$ws.Cells.Item(47, 2).Value = 601
$ws.Cells.Item(47, 3).Value = 137
$ws.Cells.Item(47, 4).Value = 384
$ws.Cells.Item(47, 5).Value = 80
$ws.Cells.Item(48, 1).Value = 'Lugo'
$ws.Cells.Item(48, 2).Value = 586
$ws.Cells.Item(48, 3).Value = 333
$ws.Cells.Item(48, 4).Value = 520
$ws.Cells.Item(48, 5).Value = 11
$ws.Cells.Item(49, 2).Value = 541
$ws.Cells.Item(49, 3).Value = 117
$ws.Cells.Item(49, 4).Value = 359
$ws.Cells.Item(49, 5).Value = 65

Write-Output "Applied provincias Spain data update"
